# v1.2.3.2 Ok - Alterar Transf Entre Contas
#
# Inserts a new row into the "Tabela1" table for a new checklist item
# ("Transf entre contas, contabilizar") right after the existing
# "Alterar Transferencia de Contas" row, and marks that existing row
# as done (Status = "Ok").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Mark the "Alterar Transferencia de Contas" row (row 19) as completed.
$ws.Range("B19").Value = "Ok"

# Insert a new physical row above row 20 so the rows below (Multiplos
# Cartoes, Front-End, etc.) shift down by one, then grow the table to
# include it.
$ws.Rows.Item(20).Insert()
$tbl.Resize($ws.Range("B3:G24"))

# Fill in the data for the newly inserted table row.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = "Transf entre contas, contabilizar"

# Keep the active cell selection pointing at the next empty status cell.
$ws.Range("B20").Select() | Out-Null
